$d = $word.ActiveDocument
$wmain = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ------------------------------------------------------------------
# Edit 1: merge the "During training..." sentence runs into one run.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "During training, I used hidden layers with 300 and 150 units respectively in model.py.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "During training, I used hidden layers with 300 and 150 units respectively in model.py.",
    2) | Out-Null

# ------------------------------------------------------------------
# Edit 2: rewrite the "Future Work" closing paragraph into a
# 5-item numbered list, keeping the _GoBack bookmark on the last item.
# ------------------------------------------------------------------
$lastP = $d.Paragraphs.Last
$r = $lastP.Range

# Insert 4 blank paragraphs ahead of the bookmark paragraph so we end
# up with 5 paragraphs total to populate.
$r.InsertParagraphBefore() | Out-Null
$r.InsertParagraphBefore() | Out-Null
$r.InsertParagraphBefore() | Out-Null
$r.InsertParagraphBefore() | Out-Null

$total = $d.Paragraphs.Count
$i1 = $total - 4
$i2 = $total - 3
$i3 = $total - 2
$i4 = $total - 1
$i5 = $total

$p1 = $d.Paragraphs.Item($i1)
$p2 = $d.Paragraphs.Item($i2)
$p3 = $d.Paragraphs.Item($i3)
$p4 = $d.Paragraphs.Item($i4)
$p5 = $d.Paragraphs.Item($i5)

# Seed paragraph 1 with the numbered-list formatting; this mints the
# numId/abstractNum definitions (numId 5) used by all five paragraphs.
$p1.Range.Text = "x"
$p1.Style = "List Paragraph"
$p1.Range.ListFormat.ApplyNumberDefault() | Out-Null

$xml1 = '<w:p ' + $wmain + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr>' +
  '<w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">From the scores over episodes plot we can see the episode score remained 0 for a long time and the average score is close to 0 until 600 episodes. This is not like a normal training process. </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/><w:r><w:t>Also</w:t></w:r><w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t xml:space="preserve"> I tried to continue the training process and see the average score went back to very low number after 900 episodes. </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/><w:r><w:t>So</w:t></w:r><w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t xml:space="preserve"> in the future work I want to modify the </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>ddpg</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> agent and model parameters to see if I can implement a steady increasing average score model which is more reasonable. </w:t></w:r>' +
  '</w:p>'
$p1.Range.InsertXML($xml1)

$xml2 = '<w:p ' + $wmain + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr>' +
  '<w:r><w:t xml:space="preserve">Modifying the hyperparameters of the model to achieve less episodes and faster training time. </w:t></w:r>' +
  '</w:p>'
$p2.Range.InsertXML($xml2)

$xml3 = '<w:p ' + $wmain + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr>' +
  '<w:r><w:t>Since I used AWS ec2 instance to train the model, I didn' + [char]8217 + 't visualize how my agent work in real environment. I will copy the weights to install it on my computer and see if there is anything I can improve from the performance of real agent.</w:t></w:r>' +
  '</w:p>'
$p3.Range.InsertXML($xml3)

$xml4 = '<w:p ' + $wmain + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr>' +
  '<w:r><w:t xml:space="preserve">I will try to complete soccer training after I feel comfortable with the current model in tennis game. </w:t></w:r>' +
  '</w:p>'
$p4.Range.InsertXML($xml4)

$xml5 = '<w:p ' + $wmain + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr>' +
  '<w:r><w:t>Compare low level model without experience replay and actor-critic with my current model to see if how these complicated algorithm affect the performance and the training speed.</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '</w:p>'
$p5.Range.InsertXML($xml5)
